$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.247862339019775
$ws.Range("B1").Value = 2.353152751922607
$ws.Range("C1").Value = 3.308939218521118
$ws.Range("D1").Value = 2.304743528366089
$ws.Range("E1").Value = 1.367507338523865
